# Update the "cryptos" price list: refreshed Price (D) / Volume(1h) (E)
# figures, plus two row swaps (TrustWalletToken<->RenderToken,
# Aave<->Aptos) from the latest GitHub Actions data pull.
#
# Numeric-looking price strings are forced to Text (NumberFormat "@")
# before assignment so Excel keeps them as literal strings (e.g. "1.000",
# "108.10") instead of silently re-typing them as numbers and dropping
# trailing zeros / precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.798.77'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '1.886.31'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.16'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4759'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2914'
$ws.Range('E8').Value = '  +6.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06581'
$ws.Range('E9').Value = '  +4.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.75'
$ws.Range('E10').Value = '  +9.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '101.36'
$ws.Range('E11').Value = '  +20.58%  '
$ws.Range('D12').Value = '1.887.97'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07580'
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.113'
$ws.Range('E14').Value = '  +3.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6579'
$ws.Range('E15').Value = '  +5.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '303.74'
$ws.Range('E16').Value = '  +31.55%  '
$ws.Range('D17').Value = '30.797.69'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.15'
$ws.Range('E18').Value = '  +4.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007574'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').Value = '2.133.93'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.115'
$ws.Range('E23').Value = '  +3.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.186'
$ws.Range('E24').Value = '  +4.86%  '
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.21'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.49'
$ws.Range('E27').Value = '  +14.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.946'
$ws.Range('E28').Value = '  +3.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1085'
$ws.Range('E29').Value = '  +6.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.353'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.163'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.978'
$ws.Range('E32').Value = '  +3.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05116'
$ws.Range('E33').Value = '  +4.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.168'
$ws.Range('E34').Value = '  +2.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7298'
$ws.Range('E35').Value = '  +3.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.715'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.706'
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9126'
$ws.Range('E39').Value = '  +4.63%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.072'
$ws.Range('E40').Value = '  +6.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '108.10'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4193'
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.631'
$ws.Range('E44').Value = '  +1.88%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.17'
$ws.Range('E45').Value = '  +7.85%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.359'
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.028'
$ws.Range('E48').Value = '  +5.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.80'
$ws.Range('E49').Value = '  +4.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05643'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('E51').Value = '  +2.28%  '
